$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 41: The Write Stuff / Enchanted Mythril Ink
$ws.Range("H41").Value = 254.05
$ws.Range("I41").Value = 71.75
$ws.Range("J41").Value = 375.58334
$ws.Range("K41").Value = 71.75
$ws.Range("L41").Value = 375.58334
$ws.Range("M41").Value = 368.25
$ws.Range("N41").Value = -1255.58334

# ALC row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 1290.125
$ws.Range("I127").Value = 748.5
$ws.Range("J127").Value = 1470.6666
$ws.Range("K127").Value = 2245.5
$ws.Range("L127").Value = 4411.9998
$ws.Range("M127").Value = 2714.5
$ws.Range("N127").Value = -14331.9998

# ALC row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 844.5867
$ws.Range("J129").Value = 883.0441
$ws.Range("L129").Value = 2649.1323
$ws.Range("N129").Value = -12649.1323

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2412.7307
$ws.Range("I138").Value = 598.4375
$ws.Range("J138").Value = 3219.0833
$ws.Range("K138").Value = 1795.3125
$ws.Range("L138").Value = 9657.249899999999
$ws.Range("M138").Value = 3344.6875
$ws.Range("N138").Value = -19937.2499

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1804.0256
$ws.Range("I32").Value = 1529.3334
$ws.Range("J32").Value = 8671.333000000001
$ws.Range("K32").Value = 1529.3334
$ws.Range("L32").Value = 8671.333000000001
$ws.Range("M32").Value = -1242.3334
$ws.Range("N32").Value = -9245.333000000001

# ARM row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 3246.8333
$ws.Range("I45").Value = 2523.3333
$ws.Range("J45").Value = 3680.9333
$ws.Range("K45").Value = 2523.3333
$ws.Range("L45").Value = 3680.9333
$ws.Range("M45").Value = -2146.3333
$ws.Range("N45").Value = -4434.933300000001

# ARM row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 3127050
$ws.Range("I63").Value = 2312.375
$ws.Range("J63").Value = 15626000
$ws.Range("K63").Value = 2312.375
$ws.Range("L63").Value = 15626000
$ws.Range("M63").Value = -1626.375
$ws.Range("N63").Value = -15627372

# ARM row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 3127050
$ws.Range("I66").Value = 2312.375
$ws.Range("J66").Value = 15626000
$ws.Range("K66").Value = 11561.875
$ws.Range("L66").Value = 78130000
$ws.Range("M66").Value = -8129.875
$ws.Range("N66").Value = -78136864

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2522
$ws.Range("I122").Value = 2456.682
$ws.Range("K122").Value = 7370.045999999999
$ws.Range("M122").Value = -4920.045999999999

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 28524.947
$ws.Range("I132").Value = 1972.7693
$ws.Range("J132").Value = 86054.664
$ws.Range("K132").Value = 5918.3079
$ws.Range("L132").Value = 258163.992
$ws.Range("M132").Value = -3388.3079
$ws.Range("N132").Value = -263223.992

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1777.9
$ws.Range("I99").Value = 1777.9
$ws.Range("K99").Value = 1777.9
$ws.Range("M99").Value = -279.9000000000001

# BSM row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 3573349.8
$ws.Range("I105").Value = 2049.6667
$ws.Range("J105").Value = 6251825
$ws.Range("K105").Value = 2049.6667
$ws.Range("L105").Value = 6251825
$ws.Range("M105").Value = -302.6667000000002
$ws.Range("N105").Value = -6255319

$ws = $wb.Worksheets.Item("CRP")
# CRP row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 2070.4583
$ws.Range("I94").Value = 1550.5
$ws.Range("J94").Value = 2330.4375
$ws.Range("K94").Value = 1550.5
$ws.Range("L94").Value = 2330.4375
$ws.Range("M94").Value = -1099.5
$ws.Range("N94").Value = -3232.4375

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 3600.2
$ws.Range("I122").Value = 4000.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 12000.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -9550.75
$ws.Range("N122").Value = -10900

# CRP row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3680.8484
$ws.Range("I134").Value = 548.93335
$ws.Range("J134").Value = 35000
$ws.Range("K134").Value = 1646.80005
$ws.Range("L134").Value = 105000
$ws.Range("M134").Value = 888.1999499999999
$ws.Range("N134").Value = -110070

$ws = $wb.Worksheets.Item("CUL")
# CUL row 63: The Next to Last Supper / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 8333.333000000001
$ws.Range("I63").Value = 7750
$ws.Range("J63").Value = 9500
$ws.Range("K63").Value = 23250
$ws.Range("L63").Value = 28500
$ws.Range("M63").Value = -22501
$ws.Range("N63").Value = -29998

# CUL row 66: Nostalgia through the Stomach (L) / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 8333.333000000001
$ws.Range("I66").Value = 7750
$ws.Range("J66").Value = 9500
$ws.Range("K66").Value = 69750
$ws.Range("L66").Value = 85500
$ws.Range("M66").Value = -66006
$ws.Range("N66").Value = -92988

# CUL row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 5027.346
$ws.Range("J68").Value = 6026.7617
$ws.Range("L68").Value = 18080.2851
$ws.Range("N68").Value = -19702.2851

# CUL row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 5027.346
$ws.Range("J71").Value = 6026.7617
$ws.Range("L71").Value = 54240.8553
$ws.Range("N71").Value = -62352.8553

# CUL row 74: The Nutcracker's Sweets / Royal Eggs
$ws.Range("H74").Value = 9850
$ws.Range("J74").Value = 9850
$ws.Range("L74").Value = 29550
$ws.Range("N74").Value = -31672

# CUL row 77: Time for a Midnight Snack (L) / Royal Eggs
$ws.Range("H77").Value = 9850
$ws.Range("J77").Value = 9850
$ws.Range("L77").Value = 88650
$ws.Range("N77").Value = -99258

# CUL row 96: Hunger Is No Game / Popoto Soba
$ws.Range("H96").Value = 503971.34
$ws.Range("J96").Value = 503971.34
$ws.Range("L96").Value = 1511914.02
$ws.Range("N96").Value = -1516032.02

# CUL row 105: Fish Box / Chirashi-zushi
$ws.Range("H105").Value = 10017.4
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 10017.4
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 30052.2
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -35294.2

# CUL row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 4799
$ws.Range("J107").Value = 1038
$ws.Range("L107").Value = 3114
$ws.Range("N107").Value = -6954

# CUL row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 799.6667
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340

# CUL row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 961.8077
$ws.Range("I122").Value = 343.55554
$ws.Range("K122").Value = 3091.99986
$ws.Range("M122").Value = -641.9998599999999

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 817.92
$ws.Range("J131").Value = 827.75256
$ws.Range("L131").Value = 2483.25768
$ws.Range("N131").Value = -12563.25768

$ws = $wb.Worksheets.Item("GSM")
# GSM row 46: Burning the Midnight Oil / Fire Brand
$ws.Range("H46").Value = 21950
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 23277.777
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 23277.777
$ws.Range("M46").Value = -9844
$ws.Range("N46").Value = -23589.777

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3862
$ws.Range("I122").Value = 2928.2222
$ws.Range("J122").Value = 6663.3335
$ws.Range("K122").Value = 8784.6666
$ws.Range("L122").Value = 19990.0005
$ws.Range("M122").Value = -6334.6666
$ws.Range("N122").Value = -24890.0005

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 7038.636
$ws.Range("I126").Value = 6719.231
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 20157.693
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -17687.693
$ws.Range("N126").Value = -27440

# GSM row 138: Orders Anonymous / White Gold Halfmask of Maiming
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1785.8572
$ws.Range("I22").Value = 2000.1666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 2000.1666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -1705.1666
$ws.Range("N22").Value = -1090

# LTW row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1785.8572
$ws.Range("I27").Value = 2000.1666
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 2000.1666
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -1893.1666
$ws.Range("N27").Value = -714

# LTW row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 2050.5557
$ws.Range("I93").Value = 2128.6667
$ws.Range("K93").Value = 2128.6667
$ws.Range("M93").Value = -880.6667000000002

# LTW row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3929.8333
$ws.Range("I122").Value = 3308.8333
$ws.Range("J122").Value = 4550.8335
$ws.Range("K122").Value = 9926.499899999999
$ws.Range("L122").Value = 13652.5005
$ws.Range("M122").Value = -7476.499899999999
$ws.Range("N122").Value = -18552.5005

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 1505.931
$ws.Range("I132").Value = 1262.2084
$ws.Range("J132").Value = 2675.8
$ws.Range("K132").Value = 3786.6252
$ws.Range("L132").Value = 8027.400000000001
$ws.Range("M132").Value = -1256.6252
$ws.Range("N132").Value = -13087.4

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 4550
$ws.Range("J96").Value = 8333.333000000001
$ws.Range("L96").Value = 8333.333000000001
$ws.Range("N96").Value = -11079.333

# WVR row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1667.0454
$ws.Range("I122").Value = 1512.8334
$ws.Range("K122").Value = 4538.5002
$ws.Range("M122").Value = -2088.5002

# WVR row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 38429
$ws.Range("J123").Value = 38429
$ws.Range("L123").Value = 38429
$ws.Range("N123").Value = -48229

# WVR row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1824.375
$ws.Range("I126").Value = 1632.5
$ws.Range("K126").Value = 4897.5
$ws.Range("M126").Value = -2427.5

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2951.7856
$ws.Range("I132").Value = 2909.6155
$ws.Range("K132").Value = 8728.8465
$ws.Range("M132").Value = -6198.8465
